$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.2720863333333334
$ws.Range("H2").Value = 0.8162590000000001
$ws.Range("I2").Value = 0.1601964374275199
$ws.Range("J2").Value = 0.2007414748833069
$ws.Range("M2").Value = 23.80409633333333
$ws.Range("N2").Value = 71.41228899999999
$ws.Range("O2").Value = 0.2035379592047788
$ws.Range("P2").Value = 0.2090753787995941
$ws.Range("Q2").Value = 6.476769289650111
$ws.Range("R2").Value = 58.290923606851
$ws.Range("S2").Value = 0.03260605594587344
$ws.Range("T2").Value = 0.0419700999020166

# Row 3
$ws.Range("G3").Value = 0.2720863333333334
$ws.Range("H3").Value = 0.8162590000000001
$ws.Range("I3").Value = 0.1601964374275199
$ws.Range("J3").Value = 0.2007414748833069
$ws.Range("O3").Value = 0.3053048738509329
$ws.Range("P3").Value = 0.3136109470643028
$ws.Range("Q3").Value = 9.715088225625667
$ws.Range("R3").Value = 87.435794030631
$ws.Range("S3").Value = 0.04890875312017783
$ws.Range("T3").Value = 0.06295472405323883

# Row 4
$ws.Range("G4").Value = 0.2720863333333334
$ws.Range("H4").Value = 0.8162590000000001
$ws.Range("I4").Value = 0.1601964374275199
$ws.Range("J4").Value = 0.2007414748833069
$ws.Range("M4").Value = 33.79564933333334
$ws.Range("N4").Value = 101.386948
$ws.Range("O4").Value = 0.2889711669362822
$ws.Range("P4").Value = 0.2968328680576918
$ws.Range("Q4").Value = 9.19533430972578
$ws.Range("R4").Value = 82.75800878753201
$ws.Range("S4").Value = 0.04629215146246554
$ws.Range("T4").Value = 0.05958666772774308

# Row 5
$ws.Range("G5").Value = 0.2720863333333334
$ws.Range("H5").Value = 0.8162590000000001
$ws.Range("I5").Value = 0.1601964374275199
$ws.Range("J5").Value = 0.2007414748833069
$ws.Range("M5").Value = 9.2924895
$ws.Range("N5").Value = 18.584979
$ws.Range("O5").Value = 0.07945583492339121
$ws.Range("P5").Value = 0.05441166469831967
$ws.Range("Q5").Value = 2.528359395593501
$ws.Range("R5").Value = 15.170156373561
$ws.Range("S5").Value = 0.01272854168755639
$ws.Range("T5").Value = 0.01092267782239665

# Row 6
$ws.Range("G6").Value = 0.2720863333333334
$ws.Range("H6").Value = 0.8162590000000001
$ws.Range("I6").Value = 0.1601964374275199
$ws.Range("J6").Value = 0.2007414748833069
$ws.Range("M6").Value = 14.353493
$ws.Range("N6").Value = 43.060479
$ws.Range("O6").Value = 0.122730165084615
$ws.Range("P6").Value = 0.1260691413800917
$ws.Range("Q6").Value = 3.905389280895667
$ws.Range("R6").Value = 35.14850352806101
$ws.Range("S6").Value = 0.01966093521144671
$ws.Range("T6").Value = 0.02530730537791175

# Row 7
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.397225
$ws.Range("H7").Value = 1.191675
$ws.Range("I7").Value = 0.2338744069853316
$ws.Range("J7").Value = 0.2930670253945926
$ws.Range("M7").Value = 23.80409633333333
$ws.Range("N7").Value = 71.41228899999999
$ws.Range("O7").Value = 0.2035379592047788
$ws.Range("P7").Value = 0.2090753787995941
$ws.Range("Q7").Value = 9.455582166008332
$ws.Range("R7").Value = 85.10023949407498
$ws.Range("S7").Value = 0.04760231950802225
$ws.Range("T7").Value = 0.0612730993480447

# Row 8
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.397225
$ws.Range("H8").Value = 1.191675
$ws.Range("I8").Value = 0.2338744069853316
$ws.Range("J8").Value = 0.2930670253945926
$ws.Range("O8").Value = 0.3053048738509329
$ws.Range("P8").Value = 0.3136109470643028
$ws.Range("Q8").Value = 14.183277319175
$ws.Range("R8").Value = 127.649495872575
$ws.Range("S8").Value = 0.07140299632161839
$ws.Range("T8").Value = 0.09190902738731625

# Row 9
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.397225
$ws.Range("H9").Value = 1.191675
$ws.Range("I9").Value = 0.2338744069853316
$ws.Range("J9").Value = 0.2930670253945926
$ws.Range("M9").Value = 33.79564933333334
$ws.Range("N9").Value = 101.386948
$ws.Range("O9").Value = 0.2889711669362822
$ws.Range("P9").Value = 0.2968328680576918
$ws.Range("Q9").Value = 13.42447680643333
$ws.Range("R9").Value = 120.8202912579
$ws.Range("S9").Value = 0.06758296030308227
$ws.Range("T9").Value = 0.0869919256810133

# Row 10
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.397225
$ws.Range("H10").Value = 1.191675
$ws.Range("I10").Value = 0.2338744069853316
$ws.Range("J10").Value = 0.2930670253945926
$ws.Range("M10").Value = 9.2924895
$ws.Range("N10").Value = 18.584979
$ws.Range("O10").Value = 0.07945583492339121
$ws.Range("P10").Value = 0.05441166469831967
$ws.Range("Q10").Value = 3.6912091416375
$ws.Range("R10").Value = 22.147254849825
$ws.Range("S10").Value = 0.01858268627423252
$ws.Range("T10").Value = 0.01594626471990451

# Row 11
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.397225
$ws.Range("H11").Value = 1.191675
$ws.Range("I11").Value = 0.2338744069853316
$ws.Range("J11").Value = 0.2930670253945926
$ws.Range("M11").Value = 14.353493
$ws.Range("N11").Value = 43.060479
$ws.Range("O11").Value = 0.122730165084615
$ws.Range("P11").Value = 0.1260691413800917
$ws.Range("Q11").Value = 5.701566256925
$ws.Range("R11").Value = 51.31409631232501
$ws.Range("S11").Value = 0.02870344457837617
$ws.Range("T11").Value = 0.03694670825831382

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.029143
$ws.Range("H12").Value = 2.058286
$ws.Range("I12").Value = 0.6059291555871485
$ws.Range("J12").Value = 0.5061914997221006
$ws.Range("M12").Value = 23.80409633333333
$ws.Range("N12").Value = 71.41228899999999
$ws.Range("O12").Value = 0.2035379592047788
$ws.Range("P12").Value = 0.2090753787995941
$ws.Range("Q12").Value = 24.49781911277566
$ws.Range("R12").Value = 146.986914676654
$ws.Range("S12").Value = 0.1233295837508831
$ws.Range("T12").Value = 0.1058321795495328

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.029143
$ws.Range("H13").Value = 2.058286
$ws.Range("I13").Value = 0.6059291555871485
$ws.Range("J13").Value = 0.5061914997221006
$ws.Range("O13").Value = 0.3053048738509329
$ws.Range("P13").Value = 0.3136109470643028
$ws.Range("Q13").Value = 36.746480131129
$ws.Range("R13").Value = 220.478880786774
$ws.Range("S13").Value = 0.1849931244091367
$ws.Range("T13").Value = 0.1587471956237478

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.029143
$ws.Range("H14").Value = 2.058286
$ws.Range("I14").Value = 0.6059291555871485
$ws.Range("J14").Value = 0.5061914997221006
$ws.Range("M14").Value = 33.79564933333334
$ws.Range("N14").Value = 101.386948
$ws.Range("O14").Value = 0.2889711669362822
$ws.Range("P14").Value = 0.2968328680576918
$ws.Range("Q14").Value = 34.78055594185467
$ws.Range("R14").Value = 208.683335651128
$ws.Range("S14").Value = 0.1750960551707344
$ws.Range("T14").Value = 0.1502542746489354

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.029143
$ws.Range("H15").Value = 2.058286
$ws.Range("I15").Value = 0.6059291555871485
$ws.Range("J15").Value = 0.5061914997221006
$ws.Range("M15").Value = 9.2924895
$ws.Range("N15").Value = 18.584979
$ws.Range("O15").Value = 0.07945583492339121
$ws.Range("P15").Value = 0.05441166469831967
$ws.Range("Q15").Value = 9.563300521498499
$ws.Range("R15").Value = 38.253202085994
$ws.Range("S15").Value = 0.0481446069616023
$ws.Range("T15").Value = 0.02754272215601852

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.029143
$ws.Range("H16").Value = 2.058286
$ws.Range("I16").Value = 0.6059291555871485
$ws.Range("J16").Value = 0.5061914997221006
$ws.Range("M16").Value = 14.353493
$ws.Range("N16").Value = 43.060479
$ws.Range("O16").Value = 0.122730165084615
$ws.Range("P16").Value = 0.1260691413800917
$ws.Range("Q16").Value = 14.771796846499
$ws.Range("R16").Value = 88.63078107899399
$ws.Range("S16").Value = 0.006450957514083397
$ws.Range("T16").Value = 0.06381512774386618
